# Update the "想去人数" (want-to-go count) figures across the four sheets
# of the 广州-漫展信息 workbook, matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 23
$ws1.Range("F3").Value  = 2750
$ws1.Range("F4").Value  = 1082
$ws1.Range("F5").Value  = 19956
$ws1.Range("F7").Value  = 2317
$ws1.Range("F11").Value = 701
$ws1.Range("F12").Value = 249
$ws1.Range("F17").Value = 272
$ws1.Range("F19").Value = 210

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 14
$ws2.Range("F7").Value  = 294
$ws2.Range("F8").Value  = 136
$ws2.Range("F16").Value = 99

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6017
$ws3.Range("F3").Value = 654

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6017
$ws4.Range("F3").Value  = 654
$ws4.Range("F6").Value  = 23
$ws4.Range("F8").Value  = 2750
$ws4.Range("F9").Value  = 1082
$ws4.Range("F10").Value = 19956
$ws4.Range("F11").Value = 14
$ws4.Range("F15").Value = 294
$ws4.Range("F16").Value = 2317
$ws4.Range("F18").Value = 136
$ws4.Range("F20").Value = 448
$ws4.Range("F21").Value = 701
$ws4.Range("F22").Value = 249
$ws4.Range("F32").Value = 272
$ws4.Range("F36").Value = 210
$ws4.Range("F37").Value = 99
$ws4.Range("F38").Value = 99

$wb.Save()
